$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21, matching the author's edit.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
